$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1333
$ws1.Range("F5").Value = 933
$ws1.Range("F6").Value = 747
$ws1.Range("G6").Value = 85
$ws1.Range("F7").Value = 215
$ws1.Range("F8").Value = 550
$ws1.Range("F9").Value = 158
$ws1.Range("F12").Value = 3085
$ws1.Range("F13").Value = 2699
$ws1.Range("F18").Value = 265
$ws1.Range("F20").Value = 5519
$ws1.Range("F22").Value = 1010
$ws1.Range("F25").Value = 411
$ws1.Range("F26").Value = 1181
$ws1.Range("F29").Value = 310

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 1167
$ws2.Range("F5").Value = 21
$ws2.Range("F13").Value = 625
$ws2.Range("F18").Value = 49
$ws2.Range("F24").Value = 287
$ws2.Range("F29").Value = 205
$ws2.Range("F30").Value = 58

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 2537
$ws3.Range("F9").Value = 1403

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 2537
$ws4.Range("F7").Value = 1403
$ws4.Range("F11").Value = 1333
$ws4.Range("F12").Value = 933
$ws4.Range("F13").Value = 747
$ws4.Range("G13").Value = 85
$ws4.Range("F14").Value = 1167
$ws4.Range("F16").Value = 215
$ws4.Range("F17").Value = 550
$ws4.Range("F18").Value = 158
$ws4.Range("F19").Value = 3085
$ws4.Range("F20").Value = 2699
$ws4.Range("F26").Value = 265
$ws4.Range("F28").Value = 5519
$ws4.Range("F30").Value = 1010
$ws4.Range("F31").Value = 625
$ws4.Range("F34").Value = 411
$ws4.Range("F40").Value = 287
$ws4.Range("F41").Value = 1181
$ws4.Range("F44").Value = 205
$ws4.Range("F45").Value = 58
$ws4.Range("F49").Value = 310
